# Add styles to the new paragraphs (see: #4)

$d = $word.ActiveDocument

# --- Create the new character styles -------------------------------------

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Datas da campaña de 2022..." run ------------

$datesText = "Datas da campaña de 2022 que usan Constelación de Hércules: 13-22 de xuño, 12-21 de xullo, 10-19 de agosto"

$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the "Grazas por participar..." run -------------

$thanksText = "Grazas por participar nesta campaña global de medida da contaminación lumínica mediante a observación das estrelas máis febles que podes albiscar. Localizando e observando a  Constelación de Hércules e comparándoa co que aparece nos mapas estelares recollidos neste documento podes saber canto contribúen á contaminación lumínica os sistemas de iluminación que hai no teu barrio ou vila. As túas achegas á base de datos en liña de GLOBE at Night (O MUNDO á Noite) servirán para documentar a calidade do ceo nocturno."

$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$found = $rng.Find.Execute($thanksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $rng.Style = "GaNParagraph"
    $rng.Collapse(0)
    $found = $rng.Find.Execute($thanksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNLinks to the "Os mapas de estrelas..." run -------------------

$mapsText = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$found = $rng.Find.Execute($mapsText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $rng.Style = "GaNLinks"
    $rng.Collapse(0)
    $found = $rng.Find.Execute($mapsText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

Write-Output "Done"
